$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as literal text (values such as
# "26.000.26" or "0.0₅ 7760" are not valid numbers, and even the ones that
# parse as numbers - e.g. "2.570" - must keep their exact text, including
# trailing zeros). Force these cells to Text format before writing so Excel
# does not silently convert the literal into a floating point number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.978.42"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.642.35"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "215.19"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "0.5066"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "0.06368"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "19.91"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "4.302"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.643.38"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "0.5474"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "0.0₅7755"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "64.35"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "26.017.05"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "4.472"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "196.94"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "9.991"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "6.164"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "1.897"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "142.76"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "0.1262"
$ws.Range("E26").Value = "  +10.33%  "
$ws.Range("D27").Value = "6.878"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "15.65"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").Value = "1.242"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "0.04906"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "3.274"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "3.208"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "1.555"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "2.379"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "0.9203"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "2.570"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "0.5555"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").Value = "1.124.06"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "0.01569"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "5.612"
$ws.Range("D42").Value = "0.8045"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "98.66"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "1.779.52"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "0.0₈118"
$ws.Range("E45").Value = "  -9.65%  "
$ws.Range("D46").Value = "0.4528"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "55.37"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "0.05189"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "7.596"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("E51").Value = "  -0.07%  "
